# Update "想去人数" (interest count) values in the "展览" and "全部类型"
# sheets to reflect refreshed scrape counts.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - column F holds the interest count.
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F6").Value  = 309
$wsExpo.Range("F7").Value  = 1112
$wsExpo.Range("F9").Value  = 6947
$wsExpo.Range("F13").Value = 7840
$wsExpo.Range("F18").Value = 2318
$wsExpo.Range("F21").Value = 272
$wsExpo.Range("F25").Value = 316
$wsExpo.Range("F28").Value = 2088
$wsExpo.Range("F30").Value = 236
$wsExpo.Range("F32").Value = 41
$wsExpo.Range("F33").Value = 545
$wsExpo.Range("F36").Value = 1416
$wsExpo.Range("F39").Value = 2151

# Sheet "全部类型" (all types) - same events, different row offsets.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F9").Value  = 309
$wsAll.Range("F10").Value = 1112
$wsAll.Range("F12").Value = 6947
$wsAll.Range("F16").Value = 7840
$wsAll.Range("F21").Value = 2318
$wsAll.Range("F24").Value = 272
$wsAll.Range("F30").Value = 316
$wsAll.Range("F33").Value = 2088
$wsAll.Range("F35").Value = 236
$wsAll.Range("F37").Value = 41
$wsAll.Range("F38").Value = 545
$wsAll.Range("F42").Value = 1416
$wsAll.Range("F45").Value = 2151
